$wb = $excel.ActiveWorkbook

# The workbook has duplicate data in sheet "展览" and sheet "全部类型".
# Update the "想去人数" (want-to-go count) values in column F for both sheets.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 138
    $ws.Range("F3").Value = 1698
    $ws.Range("F4").Value = 29
    $ws.Range("F7").Value = 155
    $ws.Range("F9").Value = 626
}
